$wb = $excel.ActiveWorkbook

# Sheet "建物" (building): the property_category column (I) was mistakenly
# left as "land" for every row. Fix it to read "building".
$wsBuilding = $wb.Worksheets.Item("建物")
$lastRow = $wsBuilding.Cells.Item($wsBuilding.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRow; $r++) {
    if ($wsBuilding.Range("I$r").Value2 -eq "land") {
        $wsBuilding.Range("I$r").Value = "building"
    }
}

# Sheet "汽車" (car): same mistake in column H -> should read "car".
$wsCar = $wb.Worksheets.Item("汽車")
$lastRowCar = $wsCar.Cells.Item($wsCar.Rows.Count, 1).End(-4162).Row
for ($r = 2; $r -le $lastRowCar; $r++) {
    if ($wsCar.Range("H$r").Value2 -eq "land") {
        $wsCar.Range("H$r").Value = "car"
    }
}

Write-Output "property_category fixed for building/car sheets"
